# Apply crypto price/volume updates to Sheet1, per commit "Updated cryptos list"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.881.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.42%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.499.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.81%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.498.90"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.51%  "
$ws.Range("E9").Value = "  -1.79%  "
$ws.Range("E10").Value = "  -0.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "8.03"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.46%  "
$ws.Range("E12").Value = "  -2.44%  "
$ws.Range("E13").Value = "  -1.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.089.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.494.36"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.837.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.65%  "
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.85%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "434.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.59%  "
$ws.Range("E23").Value = "  -3.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.635.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.82%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000120"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.33%  "
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.69%  "
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("E33").Value = "  -3.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.491.62"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.95"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.46%  "
$ws.Range("E37").Value = "  -4.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.99"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "170.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -10.48%  "
$ws.Range("E44").Value = "  -2.30%  "
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.85"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.36%  "
$ws.Range("E48").Value = "  -5.39%  "
$ws.Range("E49").Value = "  -2.77%  "
$ws.Range("E50").Value = "  -5.11%  "
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.969"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.67%  "
